$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.330.36'
$ws.Range("E2").Value = '  +6.57%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.790.35'
$ws.Range("E3").Value = '  +22.80%  '

# Row 4
$ws.Range("E4").Value = '  -0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.63'
$ws.Range("E5").Value = '  +7.74%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '179.68'
$ws.Range("E6").Value = '  +1.45%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.789.76'
$ws.Range("E7").Value = '  +22.90%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  +6.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  +11.13%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.44'
$ws.Range("E11").Value = '  +1.15%  '

# Row 12
$ws.Range("E12").Value = '  +7.76%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.58'
$ws.Range("E13").Value = '  +13.28%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000260'
$ws.Range("E14").Value = '  +8.74%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.417.44'
$ws.Range("E15").Value = '  +22.67%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.779.05'
$ws.Range("E16").Value = '  +22.46%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '71.457.60'
$ws.Range("E17").Value = '  +6.79%  '

# Row 18
$ws.Range("E18").Value = '  +1.47%  '

# Row 19
$ws.Range("E19").Value = '  +7.92%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '524.71'
$ws.Range("E20").Value = '  +8.22%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.76'
$ws.Range("E21").Value = '  +1.81%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.38'
$ws.Range("E22").Value = '  +22.45%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.747'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.52'
$ws.Range("E24").Value = '  +6.27%  '

# Row 25
$ws.Range("E25").Value = '  +11.40%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.54'
$ws.Range("E26").Value = '  +7.99%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.06'
$ws.Range("E27").Value = '  +9.44%  '

# Row 28
$ws.Range("E28").Value = '  +0.12%  '

# Row 29
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0000121'
$ws.Range("E29").Value = '  +29.56%  '

# Row 30
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.52'
$ws.Range("E30").Value = '  +10.61%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.92'
$ws.Range("E31").Value = '  +13.28%  '

# Row 32
$ws.Range("E32").Value = '  +2.72%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '32.19'
$ws.Range("E33").Value = '  +15.25%  '

# Row 34
$ws.Range("E34").Value = '  +4.00%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.08%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("E36").Value = '  +13.04%  '

# Row 37
$ws.Range("E37").Value = '  +10.70%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.24'
$ws.Range("E38").Value = '  +11.89%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.340'
$ws.Range("E39").Value = '  +10.04%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.132'
$ws.Range("E40").Value = '  +7.69%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.93'
$ws.Range("E41").Value = '  +6.08%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '431.79'
$ws.Range("E42").Value = '  +18.11%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.161.35'
$ws.Range("E43").Value = '  +13.29%  '

# Row 44
$ws.Range("E44").Value = '  +7.66%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '43.94'
$ws.Range("E45").Value = '  -6.65%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.80'
$ws.Range("E46").Value = '  +3.75%  '

# Row 47
$ws.Range("E47").Value = '  +7.11%  '

# Row 48
$ws.Range("E48").Value = '  +10.13%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '140.67'
$ws.Range("E49").Value = '  +3.93%  '

# Row 50
$ws.Range("E50").Value = '  +0.02%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.46'
$ws.Range("E51").Value = '  +7.52%  '
